# Scheduled runner update: refresh market-board derived profit figures
# (currentAveragePrice / NQ / HQ / LevePrice / LeveProfit columns, H:N)
# across the ALC/ARM/BSM/CUL/GSM/LTW/WVR Leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2926.8928
$ws.Range("I86").Value = 2553.6875
$ws.Range("J86").Value = 3424.5
$ws.Range("K86").Value = 2553.6875
$ws.Range("L86").Value = 3424.5
$ws.Range("M86").Value = -1430.6875
$ws.Range("N86").Value = -5670.5

$ws.Range("H89").Value = 2926.8928
$ws.Range("I89").Value = 2553.6875
$ws.Range("J89").Value = 3424.5
$ws.Range("K89").Value = 12768.4375
$ws.Range("L89").Value = 17122.5
$ws.Range("M89").Value = -7152.4375
$ws.Range("N89").Value = -28354.5

$ws.Range("H98").Value = 2407
$ws.Range("I98").Value = 2010.25
$ws.Range("J98").Value = 3313.8572
$ws.Range("K98").Value = 2010.25
$ws.Range("L98").Value = 3313.8572
$ws.Range("M98").Value = -512.25
$ws.Range("N98").Value = -6309.8572

$ws.Range("H116").Value = 3476.8462
$ws.Range("I116").Value = 3471.4285
$ws.Range("J116").Value = 3483.1667
$ws.Range("K116").Value = 3471.4285
$ws.Range("L116").Value = 3483.1667
$ws.Range("M116").Value = -29.42849999999999

$ws.Range("H122").Value = 2407
$ws.Range("I122").Value = 2010.25
$ws.Range("J122").Value = 3313.8572
$ws.Range("K122").Value = 6030.75
$ws.Range("L122").Value = 9941.571599999999
$ws.Range("M122").Value = -3580.75
$ws.Range("N122").Value = -14841.5716

$ws.Range("H137").Value = 1239.6471
$ws.Range("I137").Value = 1299.8334
$ws.Range("J137").Value = 1095.2
$ws.Range("K137").Value = 3899.5002
$ws.Range("L137").Value = 3285.6
$ws.Range("M137").Value = -1349.5002
$ws.Range("N137").Value = -8385.6

$ws.Range("H138").Value = 1300.73
$ws.Range("I138").Value = 665.7222
$ws.Range("J138").Value = 1657.9219
$ws.Range("K138").Value = 1997.1666
$ws.Range("L138").Value = 4973.7657
$ws.Range("M138").Value = 3142.8334
$ws.Range("N138").Value = -15253.7657

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3150.805
$ws.Range("I32").Value = 3275.7568
$ws.Range("J32").Value = 1995
$ws.Range("K32").Value = 3275.7568
$ws.Range("L32").Value = 1995
$ws.Range("M32").Value = -2988.7568
$ws.Range("N32").Value = -2569

$ws.Range("H61").Value = 2444.5715
$ws.Range("I61").Value = 2022.4
$ws.Range("J61").Value = 3500
$ws.Range("K61").Value = 2022.4
$ws.Range("L61").Value = 3500
$ws.Range("M61").Value = -1810.4
$ws.Range("N61").Value = -3924

$ws.Range("H63").Value = 2324
$ws.Range("I63").Value = 2324
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2324
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1638

$ws.Range("H66").Value = 2324
$ws.Range("I66").Value = 2324
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 11620
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -8188

$ws.Range("H107").Value = 30000
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 30000
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 30000
$ws.Range("N107").Value = -37680

$ws.Range("H110").Value = 1954.2778
$ws.Range("I110").Value = 1421.4445
$ws.Range("J110").Value = 2487.111
$ws.Range("K110").Value = 1421.4445
$ws.Range("L110").Value = 2487.111
$ws.Range("M110").Value = 623.5554999999999
$ws.Range("N110").Value = -6577.111

$ws.Range("H132").Value = 2395
$ws.Range("I132").Value = 1985.1364
$ws.Range("J132").Value = 4198.4
$ws.Range("K132").Value = 5955.4092
$ws.Range("L132").Value = 12595.2
$ws.Range("M132").Value = -3425.4092

$ws.Range("H136").Value = 2444.5715
$ws.Range("I136").Value = 2022.4
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 6067.200000000001
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -3517.200000000001
$ws.Range("N136").Value = -15600

$ws.Range("H138").Value = 125000
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 125000
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 125000
$ws.Range("N138").Value = -135280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3691.087
$ws.Range("I86").Value = 4355.4
$ws.Range("J86").Value = 2445.5
$ws.Range("K86").Value = 4355.4
$ws.Range("L86").Value = 2445.5
$ws.Range("M86").Value = -3232.4
$ws.Range("N86").Value = -4691.5

$ws.Range("H89").Value = 3691.087
$ws.Range("I89").Value = 4355.4
$ws.Range("J89").Value = 2445.5
$ws.Range("K89").Value = 21777
$ws.Range("L89").Value = 12227.5
$ws.Range("M89").Value = -16161
$ws.Range("N89").Value = -23459.5

$ws.Range("H94").Value = 35715176
$ws.Range("I94").Value = 41667536
$ws.Range("J94").Value = 1010
$ws.Range("K94").Value = 41667536
$ws.Range("L94").Value = 1010
$ws.Range("M94").Value = -41667085

$ws.Range("H107").Value = 1729.52
$ws.Range("I107").Value = 1325
$ws.Range("J107").Value = 2589.125
$ws.Range("K107").Value = 1325
$ws.Range("L107").Value = 2589.125
$ws.Range("M107").Value = 595
$ws.Range("N107").Value = -6429.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H119").Value = 7208.6
$ws.Range("I119").Value = 2282.25
$ws.Range("J119").Value = 9000
$ws.Range("K119").Value = 6846.75
$ws.Range("L119").Value = 27000
$ws.Range("M119").Value = -2008.75

$ws.Range("H122").Value = 342.16666
$ws.Range("I122").Value = 319.6
$ws.Range("J122").Value = 455
$ws.Range("K122").Value = 2876.4
$ws.Range("L122").Value = 4095
$ws.Range("M122").Value = -426.4000000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("N21").ClearContents()

$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").ClearContents()
$ws.Range("N30").ClearContents()

$ws.Range("H122").Value = 252034.5
$ws.Range("I122").Value = 1902.3334
$ws.Range("J122").Value = 502166.66
$ws.Range("K122").Value = 5707.0002
$ws.Range("L122").Value = 1506499.98
$ws.Range("M122").Value = -3257.0002
$ws.Range("N122").Value = -1511399.98

$ws.Range("H132").Value = 2297.3157
$ws.Range("I132").Value = 2045.5
$ws.Range("J132").Value = 3002.4
$ws.Range("K132").Value = 6136.5
$ws.Range("L132").Value = 9007.200000000001
$ws.Range("M132").Value = -3606.5
$ws.Range("N132").Value = -14067.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1446.6923
$ws.Range("I16").Value = 1454.5
$ws.Range("J16").Value = 1440
$ws.Range("K16").Value = 1454.5
$ws.Range("L16").Value = 1440
$ws.Range("M16").Value = -1284.5
$ws.Range("N16").Value = -1780

$ws.Range("H22").Value = 1323.6875
$ws.Range("I22").Value = 1328
$ws.Range("J22").Value = 1316.5
$ws.Range("K22").Value = 1328
$ws.Range("L22").Value = 1316.5
$ws.Range("M22").Value = -1033
$ws.Range("N22").Value = -1906.5

$ws.Range("H27").Value = 1323.6875
$ws.Range("I27").Value = 1328
$ws.Range("J27").Value = 1316.5
$ws.Range("K27").Value = 1328
$ws.Range("L27").Value = 1316.5
$ws.Range("M27").Value = -1221
$ws.Range("N27").Value = -1530.5

$ws.Range("H46").Value = 2325.75
$ws.Range("I46").Value = 1900.5
$ws.Range("J46").Value = 2751
$ws.Range("K46").Value = 1900.5
$ws.Range("L46").Value = 2751
$ws.Range("M46").Value = -1712.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 500
$ws.Range("I39").Value = 500
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 500
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -87
$ws.Range("N39").ClearContents()

$ws.Range("H132").Value = 4134.091
$ws.Range("I132").Value = 3403.5334
$ws.Range("J132").Value = 5699.5713
$ws.Range("K132").Value = 10210.6002
$ws.Range("L132").Value = 17098.7139
$ws.Range("M132").Value = -7680.600199999999
$ws.Range("N132").Value = -22158.7139

$ws.Range("H136").Value = 1053
$ws.Range("I136").Value = 836.6667
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 2510.0001
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = 39.9998999999998
$ws.Range("N136").Value = -14100
